$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 749
$wsExhibition.Range("F4").Value = 1487
$wsExhibition.Range("F6").Value = 93
$wsExhibition.Range("F7").Value = 144
$wsExhibition.Range("F8").Value = 6207
$wsExhibition.Range("F12").Value = 5140
$wsExhibition.Range("F17").Value = 360
$wsExhibition.Range("F18").Value = 64
$wsExhibition.Range("F22").Value = 3659

# 演出 (Performance) sheet
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 79

# 全部类型 (All Types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 79
$wsAll.Range("F4").Value = 749
$wsAll.Range("F5").Value = 1487
$wsAll.Range("F7").Value = 93
$wsAll.Range("F8").Value = 144
$wsAll.Range("F9").Value = 6207
$wsAll.Range("F13").Value = 5140
$wsAll.Range("F18").Value = 360
$wsAll.Range("F19").Value = 64
$wsAll.Range("F23").Value = 3659
